# Applies the "output generated at 456a3b4" update to 广州-漫展信息.xlsx
# (Guangzhou convention-info workbook): refreshed "want to go" counters
# (column F) and a handful of min-price (G) tweaks pulled from the live
# bilibili listings, plus a date correction for the "世界计划25时" meetup
# (moved from 2024-08-27 to 2024-08-24), which also re-sorts its row in
# the combined "全部类型" sheet relative to the adjacent "春日计划2024"
# concert row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F4").Value  = 193
$ws.Range("F5").Value  = 1052
$ws.Range("F7").Value  = 2709
$ws.Range("F9").Value  = 1328
$ws.Range("F11").Value = 637
$ws.Range("F16").Value = 759
$ws.Range("F17").Value = 801
$ws.Range("F19").Value = 549
$ws.Range("F20").Value = 1147
$ws.Range("F22").Value = 667
$ws.Range("F25").Value = 324
$ws.Range("F26").Value = 317
$ws.Range("F28").Value = 638
$ws.Range("F29").Value = 6883
$ws.Range("F30").Value = 505
$ws.Range("F35").Value = 1661
$ws.Range("F37").Value = 123

# Row 39 ("世界计划25时主题同人茶会×晓山瑞希生日会"): event date moved
# 2024-08-27 -> 2024-08-24. B39 is plain text ("2024-08-24" would
# otherwise auto-parse as a date), so force a text number format first.
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "2024-08-24"
$ws.Range("B39").NumberFormat = "General"
$ws.Range("E39").Value = "2024.08.24 10:00-08.24 16:30"
$ws.Range("G39").Value = 58

$ws.Range("F41").Value = 157

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")

$ws.Range("F4").Value  = 4
$ws.Range("F12").Value = 204
$ws.Range("F14").Value = 45
$ws.Range("F18").Value = 221

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")

$ws.Range("F2").Value = 758

# ---------------------------------------------------------------------
# Sheet "全部类型" (all categories combined)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F5").Value  = 1052
$ws.Range("F6").Value  = 2709
$ws.Range("F7").Value  = 1328
$ws.Range("F9").Value  = 637
$ws.Range("F13").Value = 4
$ws.Range("F16").Value = 801
$ws.Range("F18").Value = 549
$ws.Range("F19").Value = 1147
$ws.Range("F23").Value = 667
$ws.Range("F26").Value = 317
$ws.Range("F28").Value = 638
$ws.Range("F29").Value = 6884
$ws.Range("F30").Value = 204
$ws.Range("F31").Value = 505

# Rows 33/34/35 need re-sorting: the "世界计划25时" row used to sort
# after "春日计划2024" (2024-08-27 > 2024-08-24 19:30), but now that its
# date moved to 2024-08-24 10:00 it sorts *before* the 19:30 concert.
# Row 33 becomes the (retimed) 世界计划25时 row, row 34 keeps its
# 春日计划2024 content (only its "want to go" count refreshes), and row
# 35 becomes a 春日计划2024 row (what used to be row 33's content).
$ws.Range("C33").Value = "广州·世界计划25时主题同人茶会×晓山瑞希生日会"
$ws.Range("D33").Value = "黄边地铁B出口黄边美食广场1层 胡桃里音乐馆(黄边店)"
$ws.Range("E33").Value = "2024.08.24 10:00-08.24 16:30"
$ws.Range("F33").Value = 148
$ws.Range("G33").Value = 58
$ws.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=87815"
$ws.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202406/rzS5X2Ko1718735908971.png"

$ws.Range("F34").Value = 45

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "2024-08-24"
$ws.Range("B35").NumberFormat = "General"
$ws.Range("C35").Value = "广州·春日计划2024——特别二次元不插电音乐会"
$ws.Range("D35").Value = "人民北路696号 广州友谊剧院"
$ws.Range("E35").Value = "2024.08.24 19:30-08.24 21:00"
$ws.Range("F35").Value = 45
$ws.Range("G35").Value = 88
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=89964"
$ws.Range("I35").Value = "//i0.hdslb.com/bfs/openplatform/202407/lHPV2n6t1722233858047.jpeg"

$ws.Range("F39").Value = 157
